$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.737.53'
$ws.Range('E2').Value = '  +5.88%  '

$ws.Range('D3').Value = '2.499.96'
$ws.Range('E3').Value = '  +2.97%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.02%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.41'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.91%  '

$ws.Range('E7').Value = '  +1.54%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.14%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.541'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.82%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.01'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.44%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0814'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.52%  '

$ws.Range('E12').Value = '  +0.81%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.34'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.84%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.16'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.21%  '

$ws.Range('D15').Value = '2.887.97'
$ws.Range('E15').Value = '  +2.96%  '

$ws.Range('D16').Value = '2.500.93'
$ws.Range('E16').Value = '  +3.78%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.845'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.65%  '

$ws.Range('D18').Value = '47.582.46'
$ws.Range('E18').Value = '  +5.74%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.68'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.48%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.58'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.59%  '

$ws.Range('D21').Value = '0.0₃0937'
$ws.Range('E21').Value = '  +1.60%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.76'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.81%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '251.41'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.94%  '

$ws.Range('E24').Value = '  +5.65%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.68%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.23'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.50%  '

$ws.Range('E27').Value = '  -0.06%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.03'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.94%  '

$ws.Range('E29').Value = '  +6.48%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.28%  '

$ws.Range('E31').Value = '  +7.67%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.49'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.55%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.04'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.78%  '

$ws.Range('E34').Value = '  +3.01%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0781'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.95%  '

$ws.Range('E36').Value = '  +0.19%  '

$ws.Range('E37').Value = '  +3.83%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.62'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.53%  '

$ws.Range('E39').Value = '  +4.30%  '

$ws.Range('E40').Value = '  +2.28%  '

$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '121.52'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.47%  '

$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.79%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.07'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.43%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0297'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.84%  '

$ws.Range('D45').Value = '1.966.47'
$ws.Range('E45').Value = '  +1.72%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.99'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.91%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.09'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.47%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.83'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.11%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.21'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.51%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.35'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +13.85%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.27'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.46%  '
